# TestRules.xlsx sheet restructuring: insert a new "RuleName" column at the
# left edge of the rule table (shifting the existing grid one column to the
# right) and append a second rule row ("Rule 2").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two header cells are merged in the original layout; unmerge first so
# each underlying cell can be written to independently while we reshuffle.
$ws.Cells.UnMerge()

# --- Drop the old column-A labels that have no home after the shift ---
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("A9").ClearContents()

# --- Re-key the rule-table header/body onto the shifted grid (col+1) ---
$ws.Range("C4").Value = 'RuleTable'

$ws.Range("C5").Value = 'CONDITON'
$ws.Range("D5").Value = 'CONDITON'
$ws.Range("E5").Value = 'CONDITON'
$ws.Range("F5").Value = 'CONDITON'
$ws.Range("G5").Value = 'CONDITION'
$ws.Range("H5").Value = 'ACTION'

$ws.Range("C6").Value = '$application : MortgageRequest'
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 'Loan from $loan'

$ws.Range("C7").Value = '$address: mailingAddress,$loan: loan,annualIncome >= $param'
$ws.Range("D7").Value = 'monthlyDebt <= $param'
$ws.Range("E7").Value = 'type == "$param"'
$ws.Range("F7").Value = 'Amount > $param'
$ws.Range("G7").Value = 'Address(state in ($param)) from $address'
$ws.Range("H7").Value = 'modify($application) { setStatus("$param")};'

$ws.Range("B8").Value = 'RuleName'
$ws.Range("C8").Value = 'Annual Income'
$ws.Range("D8").Value = 'Monthly Debt'
$ws.Range("E8").Value = 'Loan type'
$ws.Range("F8").Value = 'Amount'
$ws.Range("G8").Value = 'State'
$ws.Range("H8").Value = 'Status'

$ws.Range("B9").Value = 'Rule 1'
$ws.Range("C9").Value = 50000
$ws.Range("D9").Value = 200
$ws.Range("E9").Value = 'Mortgage'
$ws.Range("F9").Value = 30000
$ws.Range("G9").Value = '“GA”,”FL”'
$ws.Range("H9").Value = 'Approved'

# --- Brand-new "Rule 2" row appended below the existing rule ---
$ws.Range("B10").Value = 'Rule 2'
$ws.Range("C10").Value = 70000
$ws.Range("D10").Value = 400
$ws.Range("E10").Value = 'Mortgage'
$ws.Range("F10").Value = 20000
$ws.Range("G10").Value = '“NY”,”CA”,”TN”'
$ws.Range("H10").Value = 'Denied'

# --- Fix up the handful of cells whose emphasis/centering moved with them ---
$ws.Range("B8").Font.Bold = $false        # RuleName header is plain, not bold
$ws.Range("H8").Font.Bold = $true         # Status header matches the other bold headers

$ws.Range("F6").HorizontalAlignment = -4108   # xlCenter, matching the merged header row
$ws.Range("F6").VerticalAlignment = -4108     # xlCenter

# --- Re-create the two merged header cells at their new location ---
$ws.Range("C6:D6").Merge() | Out-Null
$ws.Range("E6:F6").Merge() | Out-Null

# --- New column H (the ACTION column) needs its own width ---
$ws.Columns.Item(8).ColumnWidth = 36.22

# --- Leave the selection where the author left it after adding the new row ---
$ws.Range("H11").Select() | Out-Null

Write-Output "TestRules sheet restructured: RuleName column inserted, Rule 2 row appended."